$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 6655
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 34
$ws.Range("G2").Value = -87
$ws.Range("H2").Value = -134
$ws.Range("I2").Value = -127
$ws.Range("J2").Value = -7
$ws.Range("K2").Value = 8559
$ws.Range("L2").Value = 5672
$ws.Range("M2").Value = 2886
$ws.Range("N2").Value = 2855
$ws.Range("O2").Value = 31
$ws.Range("P2").Value = 120
$ws.Range("Q2").Value = 483
$ws.Range("R2").Value = -556
$ws.Range("S2").Value = -108
$ws.Range("T2").Value = 362
$ws.Range("U2").Value = 121
$ws.Range("V2").Value = 4202
$ws.Range("W2").Value = 0.51
$ws.Range("X2").Value = -2.01
$ws.Range("Y2").Value = -4.28
$ws.Range("Z2").Value = -1.55
$ws.Range("AA2").Value = 196.53
$ws.Range("AB2").Value = 2302.07
$ws.Range("AC2").Value = -528
$ws.Range("AD2").Value = -15.39
$ws.Range("AE2").Value = 11897
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 60
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = -11.37
$ws.Range("AJ2").Value = 24000000

$ws.Range("D3").Value = 6544
$ws.Range("E3").Value = 494
$ws.Range("F3").Value = 494
$ws.Range("G3").Value = 448
$ws.Range("H3").Value = 286
$ws.Range("I3").Value = 298
$ws.Range("J3").Value = -11
$ws.Range("K3").Value = 8464
$ws.Range("L3").Value = 5300
$ws.Range("M3").Value = 3164
$ws.Range("N3").Value = 3145
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 120
$ws.Range("Q3").Value = 855
$ws.Range("R3").Value = 61
$ws.Range("S3").Value = -728
$ws.Range("T3").Value = 247
$ws.Range("U3").Value = 608
$ws.Range("V3").Value = 3623
$ws.Range("W3").Value = 7.55
$ws.Range("X3").Value = 4.38
$ws.Range("Y3").Value = 9.93
$ws.Range("Z3").Value = 3.36
$ws.Range("AA3").Value = 167.49
$ws.Range("AB3").Value = 2509.27
$ws.Range("AC3").Value = 1241
$ws.Range("AD3").Value = 11.72
$ws.Range("AE3").Value = 13103
$ws.Range("AF3").Value = 1.11
$ws.Range("AG3").Value = 70
$ws.Range("AH3").Value = 0.48
$ws.Range("AI3").Value = 5.64
$ws.Range("AJ3").Value = 24000000

$ws.Range("D4").Value = 6943
$ws.Range("E4").Value = 764
$ws.Range("F4").Value = 764
$ws.Range("G4").Value = 636
$ws.Range("H4").Value = 422
$ws.Range("I4").Value = 432
$ws.Range("J4").Value = -10
$ws.Range("K4").Value = 8295
$ws.Range("L4").Value = 4753
$ws.Range("M4").Value = 3542
$ws.Range("N4").Value = 3564
$ws.Range("O4").Value = -22
$ws.Range("P4").Value = 120
$ws.Range("Q4").Value = 1008
$ws.Range("R4").Value = -360
$ws.Range("S4").Value = -847
$ws.Range("T4").Value = 233
$ws.Range("U4").Value = 775
$ws.Range("V4").Value = 3011
$ws.Range("W4").Value = 11
$ws.Range("X4").Value = 6.08
$ws.Range("Y4").Value = 12.89
$ws.Range("Z4").Value = 5.04
$ws.Range("AA4").Value = 134.2
$ws.Range("AB4").Value = 2848.53
$ws.Range("AC4").Value = 1802
$ws.Range("AD4").Value = 9.49
$ws.Range("AE4").Value = 14851
$ws.Range("AF4").Value = 1.15
$ws.Range("AG4").Value = 140
$ws.Range("AH4").Value = 0.82
$ws.Range("AI4").Value = 7.77
$ws.Range("AJ4").Value = 24000000

$ws.Range("D5").Value = 7249
$ws.Range("E5").Value = 545
$ws.Range("F5").Value = 545
$ws.Range("G5").Value = 442
$ws.Range("H5").Value = 347
$ws.Range("I5").Value = 354
$ws.Range("J5").Value = -7
$ws.Range("K5").Value = 8651
$ws.Range("L5").Value = 4870
$ws.Range("M5").Value = 3781
$ws.Range("N5").Value = 3809
$ws.Range("O5").Value = -27
$ws.Range("P5").Value = 120
$ws.Range("Q5").Value = 551
$ws.Range("R5").Value = -497
$ws.Range("S5").Value = 47
$ws.Range("T5").Value = 435
$ws.Range("U5").Value = 116
$ws.Range("V5").Value = 3075
$ws.Range("W5").Value = 7.51
$ws.Range("X5").Value = 4.79
$ws.Range("Y5").Value = 9.6
$ws.Range("Z5").Value = 4.1
$ws.Range("AA5").Value = 128.79
$ws.Range("AB5").Value = 3114.17
$ws.Range("AC5").Value = 1475
$ws.Range("AD5").Value = 18.31
$ws.Range("AE5").Value = 15869
$ws.Range("AF5").Value = 1.7
$ws.Range("AG5").Value = 140
$ws.Range("AH5").Value = 0.52
$ws.Range("AI5").Value = 9.49
$ws.Range("AJ5").Value = 24000000

$ws.Range("D6").Value = 7887
$ws.Range("E6").Value = 543
$ws.Range("F6").Value = 543
$ws.Range("G6").Value = 463
$ws.Range("H6").Value = 378
$ws.Range("I6").Value = 381
$ws.Range("K6").Value = 8725
$ws.Range("L6").Value = 4613
$ws.Range("M6").Value = 4112
$ws.Range("N6").Value = 4119
$ws.Range("P6").Value = 120
$ws.Range("Q6").Value = 582
$ws.Range("R6").Value = -148
$ws.Range("S6").Value = -379
$ws.Range("T6").Value = 289
$ws.Range("U6").Value = 293
$ws.Range("V6").Value = 2851
$ws.Range("W6").Value = 6.88
$ws.Range("X6").Value = 4.79
$ws.Range("Y6").Value = 9.609999999999999
$ws.Range("Z6").Value = 4.35
$ws.Range("AA6").Value = 112.2
$ws.Range("AB6").Value = 3373.82
$ws.Range("AC6").Value = 1587
$ws.Range("AD6").Value = 12.01
$ws.Range("AE6").Value = 17163
$ws.Range("AF6").Value = 1.11
$ws.Range("AG6").Value = 160
$ws.Range("AH6").Value = 0.84
$ws.Range("AI6").Value = 10.08
$ws.Range("AJ6").Value = 24000000

$ws.Range("D7").Value = 8160
$ws.Range("E7").Value = 637
$ws.Range("G7").Value = 498
$ws.Range("H7").Value = 343
$ws.Range("I7").Value = 333
$ws.Range("K7").Value = 9615
$ws.Range("L7").Value = 5157
$ws.Range("M7").Value = 4458
$ws.Range("N7").Value = 4420
$ws.Range("P7").Value = 120
$ws.Range("Q7").Value = 464
$ws.Range("R7").Value = -218
$ws.Range("S7").Value = 14
$ws.Range("T7").Value = 201
$ws.Range("U7").Value = 262
$ws.Range("W7").Value = 7.8
$ws.Range("X7").Value = 4.2
$ws.Range("Y7").Value = 7.8
$ws.Range("Z7").Value = 3.74
$ws.Range("AA7").Value = 115.67
$ws.Range("AC7").Value = 1388
$ws.Range("AD7").Value = 10.34
$ws.Range("AE7").Value = 18417
$ws.Range("AF7").Value = 0.78
$ws.Range("AG7").Value = 153
$ws.Range("AH7").Value = 1.07
$ws.Range("AI7").Value = 11.05

$ws.Range("D8").Value = 9062
$ws.Range("E8").Value = 823
$ws.Range("G8").Value = 761
$ws.Range("H8").Value = 543
$ws.Range("I8").Value = 530
$ws.Range("K8").Value = 10078
$ws.Range("L8").Value = 5038
$ws.Range("M8").Value = 5040
$ws.Range("N8").Value = 4980
$ws.Range("P8").Value = 120
$ws.Range("Q8").Value = 568
$ws.Range("R8").Value = -188
$ws.Range("S8").Value = -134
$ws.Range("T8").Value = 161
$ws.Range("U8").Value = 422
$ws.Range("W8").Value = 9.09
$ws.Range("X8").Value = 6
$ws.Range("Y8").Value = 11.28
$ws.Range("Z8").Value = 5.59
$ws.Range("AA8").Value = 99.95
$ws.Range("AC8").Value = 2208
$ws.Range("AD8").Value = 5.95
$ws.Range("AE8").Value = 20750
$ws.Range("AF8").Value = 0.63
$ws.Range("AG8").Value = 150
$ws.Range("AH8").Value = 1.14
$ws.Range("AI8").Value = 6.79

$ws.Range("D9").Value = 9982
$ws.Range("E9").Value = 1056
$ws.Range("G9").Value = 982
$ws.Range("H9").Value = 744
$ws.Range("I9").Value = 724
$ws.Range("K9").Value = 10780
$ws.Range("L9").Value = 5034
$ws.Range("M9").Value = 5746
$ws.Range("N9").Value = 5800
$ws.Range("P9").Value = 120
$ws.Range("Q9").Value = 944
$ws.Range("R9").Value = -206
$ws.Range("S9").Value = -159
$ws.Range("T9").Value = 173
$ws.Range("U9").Value = 528
$ws.Range("W9").Value = 10.57
$ws.Range("X9").Value = 7.46
$ws.Range("Y9").Value = 13.44
$ws.Range("Z9").Value = 7.14
$ws.Range("AA9").Value = 87.62
$ws.Range("AC9").Value = 3019
$ws.Range("AD9").Value = 4.36
$ws.Range("AE9").Value = 24167
$ws.Range("AF9").Value = 0.54
$ws.Range("AG9").Value = 150
$ws.Range("AH9").Value = 1.14
$ws.Range("AI9").Value = 4.97
